# mie_raw-METADATA.xlsx -- "prep for Est & Coasts" edit
#
# - Data collection method text expanded to clarify the observation is of
#   above-ground vegetation.
# - "Values (range, description)" for the `estuary` attribute row was filled
#   in with the estuary code legend (ere / lqre / nre).
# - Selection / scroll position left where the author last clicked (I13),
#   and the previously-blank estuary row grew to fit its new text, nudging
#   the row heights of the rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- content edits -------------------------------------------------------

$ws.Range("B8").Value = "visual observation of species presence and abundance in above-ground vegetation"
$ws.Range("F13").Value = "ere (Englishman River Estuary, not included in analysis), lqre (Little Qualicum River Estuary), nre (Nanaimo River Estuary)"

# --- row heights (table rows reflow once the new text wraps) -------------

$ws.Rows.Item(13).RowHeight = 24.6
$ws.Rows.Item(14).RowHeight = 36.6
$ws.Rows.Item(15).RowHeight = 96.6
$ws.Rows.Item(16).RowHeight = 51
$ws.Rows.Item(17).RowHeight = 48.6
$ws.Rows.Item(18).RowHeight = 84.6

# --- view state: scroll back to top, select I13 ---------------------------

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I13").Select() | Out-Null
